# Refined metadata to be additional tab
#
# 1. Update the "time_taken" timestamps (column F) on the "data" sheet.
# 2. Add a new "metadata" worksheet, placed after "data", with a header row
#    and one data row describing the panel query.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- 1. Refresh the time_taken column on the data sheet ---------------------
$timestamps = @(
    "2021-10-05 14:19:13.106474",
    "2021-10-05 14:19:13.106482",
    "2021-10-05 14:19:13.106485",
    "2021-10-05 14:19:13.106488",
    "2021-10-05 14:19:13.106491",
    "2021-10-05 14:19:13.106494",
    "2021-10-05 14:19:13.106496",
    "2021-10-05 14:19:13.106499",
    "2021-10-05 14:19:13.106502",
    "2021-10-05 14:19:13.106504",
    "2021-10-05 14:19:13.106507",
    "2021-10-05 14:19:13.106509",
    "2021-10-05 14:19:13.106512",
    "2021-10-05 14:19:13.106514",
    "2021-10-05 14:19:13.106517"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Range("F$row").Value = $timestamps[$i]
}

# --- 2. Add the new "metadata" sheet, right after "data" --------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the bold/bordered header style already used by the "data" sheet
# (row 1, columns B:F) instead of constructing a brand-new style.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Atypical haemolytic uraemic syndrome"
$metaSheet.Range("C2").Value = 139
# Enter "2.9" as a formula literal (forces text, not a number) then collapse
# it down to a plain value, so the stored number format stays untouched.
$metaSheet.Range("D2").Formula = '="2.9"'
$metaSheet.Range("D2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$metaSheet.Range("E2").Value = "2021-07-13T10:33:32.063214Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:13.103146"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/139/?format=json"
